$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reorders (mirrors the shared-string table reorder in the source diff) ---
# Tailandia / Polonia / Chile / Rumania block: Polonia now precedes Chile
$ws.Range("A34").Value2 = "Polonia"
$ws.Range("A35").Value2 = "Chile"

# Vietnam / Islas Feroe / Moldavia / Chipre / Albania / Tunez block: Islas Feroe moved up
$ws.Range("A84").Value2 = "Islas Feroe"
$ws.Range("A85").Value2 = "Moldavia"
$ws.Range("A86").Value2 = "Republica de Chipre"
$ws.Range("A87").Value2 = "Albania"

# Cuba / Mauricio / Nigeria / Consejo Danes / Kirguistan / Ruanda block: Mauricio moved up
$ws.Range("A115").Value2 = "Mauricio"
$ws.Range("A116").Value2 = "Nigeria"
$ws.Range("A117").Value2 = "Consejo Danes para los Refugiados"
$ws.Range("A118").Value2 = "Kirguistan"

# --- Updated statistics (B Casos totales, C Nuevos casos, D Casos activos, E Recuperados, F Casos criticos, G Muertes hoy, H Muertes) ---

$ws.Range("B11").Value2 = 9991
$ws.Range("C11").Value2 = 114
$ws.Range("D11").Value2 = 131
$ws.Range("E11").Value2 = 9727
$ws.Range("F11").Value2 = 141
$ws.Range("G11").Value2 = 11
$ws.Range("H11").Value2 = 133

$ws.Range("B15").Value2 = 5394
$ws.Range("C15").Value2 = 111
$ws.Range("D15").Value2 = 9
$ws.Range("E15").Value2 = 5355
$ws.Range("F15").Value2 = 26
$ws.Range("G15").Value2 = 2
$ws.Range("H15").Value2 = 30

$ws.Range("B17").Value2 = 2868
$ws.Range("C17").Value2 = 2
$ws.Range("D17").Value2 = 6
$ws.Range("E17").Value2 = 2849
$ws.Range("F17").Value2 = 44
$ws.Range("G17").Value2 = 1
$ws.Range("H17").Value2 = 13

$ws.Range("B27").Value2 = 1497
$ws.Range("C27").Value2 = 103
$ws.Range("D27").Value2 = 10
$ws.Range("E27").Value2 = 1484
$ws.Range("F27").Value2 = 19
$ws.Range("G27").Value2 = 0
$ws.Range("H27").Value2 = 3

$ws.Range("B34").Value2 = 927
$ws.Range("C34").Value2 = 26
$ws.Range("D34").Value2 = 1
$ws.Range("E34").Value2 = 914
$ws.Range("F34").Value2 = 3
$ws.Range("G34").Value2 = 2
$ws.Range("H34").Value2 = 12

$ws.Range("B35").Value2 = 922
$ws.Range("C35").Value2 = 0
$ws.Range("D35").Value2 = 17
$ws.Range("E35").Value2 = 903
$ws.Range("F35").Value2 = 7
$ws.Range("G35").Value2 = 0
$ws.Range("H35").Value2 = 2

$ws.Range("B36").Value2 = 794
$ws.Range("C36").Value2 = 0
$ws.Range("D36").Value2 = 79
$ws.Range("E36").Value2 = 702
$ws.Range("F36").Value2 = 15
$ws.Range("G36").Value2 = 1
$ws.Range("H36").Value2 = 13

$ws.Range("B45").Value2 = 558
$ws.Range("C45").Value2 = 0
$ws.Range("D45").Value2 = 156
$ws.Range("E45").Value2 = 400
$ws.Range("F45").Value2 = 17
$ws.Range("G45").Value2 = 0
$ws.Range("H45").Value2 = 2

$ws.Range("B70").Value2 = 220
$ws.Range("C70").Value2 = 2
$ws.Range("D70").Value2 = 4
$ws.Range("E70").Value2 = 213
$ws.Range("F70").Value2 = 8
$ws.Range("G70").Value2 = 0
$ws.Range("H70").Value2 = 3

$ws.Range("B84").Value2 = 132
$ws.Range("C84").Value2 = 10
$ws.Range("D84").Value2 = 38
$ws.Range("E84").Value2 = 94
$ws.Range("F84").Value2 = 2
$ws.Range("G84").Value2 = 0
$ws.Range("H84").Value2 = 0

$ws.Range("B85").Value2 = 125
$ws.Range("C85").Value2 = 0
$ws.Range("D85").Value2 = 2
$ws.Range("E85").Value2 = 122
$ws.Range("F85").Value2 = 10
$ws.Range("G85").Value2 = 0
$ws.Range("H85").Value2 = 1

$ws.Range("B86").Value2 = 124
$ws.Range("C86").Value2 = 0
$ws.Range("D86").Value2 = 3
$ws.Range("E86").Value2 = 118
$ws.Range("F86").Value2 = 3
$ws.Range("G86").Value2 = 0
$ws.Range("H86").Value2 = 3

$ws.Range("B87").Value2 = 123
$ws.Range("C87").Value2 = 0
$ws.Range("D87").Value2 = 10
$ws.Range("E87").Value2 = 108
$ws.Range("F87").Value2 = 2
$ws.Range("G87").Value2 = 0
$ws.Range("H87").Value2 = 5

$ws.Range("B93").Value2 = 102
$ws.Range("C93").Value2 = 0
$ws.Range("D93").Value2 = 3
$ws.Range("E93").Value2 = 99
$ws.Range("F93").Value2 = 2
$ws.Range("G93").Value2 = 0
$ws.Range("H93").Value2 = 0

$ws.Range("B115").Value2 = 48
$ws.Range("C115").Value2 = 6
$ws.Range("D115").Value2 = 0
$ws.Range("E115").Value2 = 46
$ws.Range("F115").Value2 = 1
$ws.Range("G115").Value2 = 0
$ws.Range("H115").Value2 = 2

$ws.Range("B116").Value2 = 46
$ws.Range("C116").Value2 = 2
$ws.Range("D116").Value2 = 2
$ws.Range("E116").Value2 = 43
$ws.Range("F116").Value2 = 0
$ws.Range("G116").Value2 = 0
$ws.Range("H116").Value2 = 1

$ws.Range("B117").Value2 = 45
$ws.Range("C117").Value2 = 0
$ws.Range("D117").Value2 = 0
$ws.Range("E117").Value2 = 43
$ws.Range("F117").Value2 = 0
$ws.Range("G117").Value2 = 0
$ws.Range("H117").Value2 = 2

$ws.Range("B118").Value2 = 42
$ws.Range("C118").Value2 = 0
$ws.Range("D118").Value2 = 0
$ws.Range("E118").Value2 = 42
$ws.Range("F118").Value2 = 0
$ws.Range("G118").Value2 = 0
$ws.Range("H118").Value2 = 0

# --- Timestamp update ---
$ws.Range("A1").Value2 = "Datos actualizados a 25 de Marzo de 2020 a las 09:16"
